$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DELIVER")

# --- Section: "CHART" block (rows 27-36 before edit) ---
# A new row is inserted right after "IMAGE3" (row 28), which becomes the new
# "IMAGE4" label row (29). The paragraph that used to sit there ("In terms of
# backlog...") is pushed down to row 30, keeping its neighbouring B:F
# (table-like) formatting.
$ws.Rows("29:29").Insert()
$ws.Range("A29").Value = "IMAGE4"
$ws.Rows("29:29").RowHeight = 16

# A second new row is inserted right after "ANUAL CHART" (now row 32) to hold
# a brand new "IMAGE5" label; this pushes the rest of the section (the
# "Based on the chart..." paragraph, the blank spacer, "GREETINGS" and the
# closing thank-you paragraph) down by one row each.
$ws.Rows("33:33").Insert()
$ws.Range("B33").Clear()
$ws.Range("A16").Copy()
$ws.Range("A33").PasteSpecial(-4122)
$ws.Range("A33").Value = "IMAGE5"

# --- Selection bookkeeping (matches the author's last-saved cursor) ---
$ws.Range("C20").Select()
